# Weekly update: prepend a new Brocoli price observation (2 rows, Primera/Segunda)
# to the Mercado Mayorista Lo Valledor de Santiago dataset.
#
# The dataset rows 1125..1199 hold the historical (date-descending "insert-at-top")
# series. New data is inserted at rows 1125-1126 and all subsequent rows are pushed
# down by 2 positions; the 2 oldest rows that fall off the bottom of the original
# range are appended as new rows 1200-1201.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary per data row.
$varCols = 4,9,10,11,12,13,15,16   # D,I,J,K,L,M,O,P

# 1) Shift existing rows 1125..1199 down by two rows, writing from the bottom
#    upward (row 1201 down to row 1127) so that source rows are always read
#    before they themselves get overwritten.
for ($r = 1201; $r -ge 1127; $r--) {
    $srcRow = $r - 2
    foreach ($c in $varCols) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($srcRow, $c).Value2
    }
}

# 2) Rows 1200 and 1201 are brand-new rows; populate the constant columns that
#    were previously implicit (copied automatically when only modifying
#    existing cells in step 1, but these two rows did not exist before).
$constCols = @{
    1  = 6
    2  = "Mercado Mayorista Lo Valledor de Santiago"
    3  = "Metropolitana"
    5  = 13
    6  = 100112023
    7  = "Brócoli"
    8  = "Sin especificar"
    14 = "`$/unidad"
    17 = 1
    18 = "Hortaliza"
}
foreach ($r in 1200,1201) {
    foreach ($c in $constCols.Keys) {
        $ws.Cells.Item($r, $c).Value = $constCols[$c]
    }
    # Match the date number format used by the rest of column D.
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

# 3) Write the new observation into rows 1125 and 1126 (the "Primera" and
#    "Segunda" quality rows for the newly reported week).
$ws.Cells.Item(1125, 4).Value = 45021
$ws.Cells.Item(1125, 9).Value = "Primera"
$ws.Cells.Item(1125, 10).Value = 11000
$ws.Cells.Item(1125, 11).Value = 500
$ws.Cells.Item(1125, 12).Value = 600
$ws.Cells.Item(1125, 13).Value = 559
$ws.Cells.Item(1125, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1125, 16).Value = 559

$ws.Cells.Item(1126, 4).Value = 45021
$ws.Cells.Item(1126, 9).Value = "Segunda"
$ws.Cells.Item(1126, 10).Value = 4500
$ws.Cells.Item(1126, 11).Value = 400
$ws.Cells.Item(1126, 12).Value = 400
$ws.Cells.Item(1126, 13).Value = 400
$ws.Cells.Item(1126, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1126, 16).Value = 400
